$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(5, 6).Value = 142
$ws.Cells.Item(6, 6).Value = 3872
$ws.Cells.Item(7, 6).Value = 514
$ws.Cells.Item(9, 6).Value = 1285
$ws.Cells.Item(10, 6).Value = 645
$ws.Cells.Item(11, 6).Value = 361
$ws.Cells.Item(13, 6).Value = 2138
$ws.Cells.Item(15, 6).Value = 638970
$ws.Cells.Item(16, 6).Value = 1571
$ws.Cells.Item(18, 6).Value = 1384
$ws.Cells.Item(21, 6).Value = 1230
$ws.Cells.Item(22, 6).Value = 2117
$ws.Cells.Item(25, 6).Value = 1508
$ws.Cells.Item(26, 6).Value = 728
$ws.Cells.Item(27, 6).Value = 1479
$ws.Cells.Item(30, 6).Value = 1062
$ws.Cells.Item(36, 6).Value = 1290
$ws.Cells.Item(46, 6).Value = 3049
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(11, 6).Value = 144128
$ws.Cells.Item(12, 6).Value = 144128
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(6, 6).Value = 225
$ws.Cells.Item(8, 6).Value = 803
$ws.Cells.Item(9, 6).Value = 1114
$ws.Cells.Item(10, 6).Value = 615
$ws.Cells.Item(11, 6).Value = 1552
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 803
$ws.Cells.Item(3, 6).Value = 615
$ws.Cells.Item(8, 6).Value = 3872
$ws.Cells.Item(10, 6).Value = 514
$ws.Cells.Item(11, 6).Value = 1285
$ws.Cells.Item(12, 6).Value = 645
$ws.Cells.Item(13, 6).Value = 361
$ws.Cells.Item(14, 6).Value = 2138
$ws.Cells.Item(16, 6).Value = 638971
$ws.Cells.Item(19, 6).Value = 1571
$ws.Cells.Item(20, 6).Value = 144128
$ws.Cells.Item(22, 6).Value = 1384
$ws.Cells.Item(25, 6).Value = 1230
$ws.Cells.Item(26, 6).Value = 2117
$ws.Cells.Item(29, 6).Value = 1508
$ws.Cells.Item(30, 6).Value = 728
$ws.Cells.Item(32, 6).Value = 1479
$ws.Cells.Item(36, 6).Value = 1062
$ws.Cells.Item(41, 6).Value = 1290
$ws.Cells.Item(51, 6).Value = 3049
